# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-locale sheets now that the
# translated files have been handed back, flips the Status from
# "Ready for handoff" to "Handed back: in sync with en-US", links the new
# target-file cell back to the source markdown doc (same as column A), and
# widens a few columns so the longer values are readable.

$wb = $excel.ActiveWorkbook

$mdName   = "5a665907-b599-4600-b7a5-34d6ec84d56c.md"
$srcUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/923f8e09e8a508c0f0c01b23f9e2921730ae76cb/e2e/$mdName"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: Status columns for both locales move to "handed back"
# (shared text cell, so both zh-cn (E2) and de-de (F2) columns change),
# and those two columns get wider to comfortably fit the longer text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.1667
$overview.Columns.Item(6).ColumnWidth = 29.1667

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("I2").Value = $mdName
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $srcUrl, $null, $null, $mdName) | Out-Null
$zhcn.Range("J2").Value = "5a665907-b599-4600-b7a5-34d6ec84d56c.0a58aa3daa72a7d1c92a129f5dfd7358dc5318fb.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-19 08:56:29"
$zhcn.Columns.Item(3).ColumnWidth = 29.1667
$zhcn.Columns.Item(9).ColumnWidth = 39.1667
$zhcn.Columns.Item(10).ColumnWidth = 39.1667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("I2").Value = $mdName
$dede.Hyperlinks.Add($dede.Range("I2"), $srcUrl, $null, $null, $mdName) | Out-Null
$dede.Range("J2").Value = "5a665907-b599-4600-b7a5-34d6ec84d56c.0a58aa3daa72a7d1c92a129f5dfd7358dc5318fb.de-de.xlf"
$dede.Range("K2").Value = "2016-08-19 08:56:36"
$dede.Columns.Item(3).ColumnWidth = 29.1667
$dede.Columns.Item(9).ColumnWidth = 39.1667
$dede.Columns.Item(10).ColumnWidth = 39.1667
